# Edit script: apply the two substantive changes described by the diff:
#  1. Slide 5's table (graphicFrame, Shapes.Item(2)) switches its table
#     style (tableStyleId) from {B3EA9435-33DD-4F00-AAAC-6DD6AF0382FD} to
#     {202BF87D-3801-40E1-A598-9F15E327B125}.
#  2. The deck's theme (ppt/theme/theme1.xml, the one actually in effect
#     for the slide master / all slides) is swapped from the "Integral"
#     (Red Violet) color palette to the stock "Office Theme" (Office)
#     color palette. (fontScheme/fmtScheme are identical between the two
#     themes in this deck, so only the 12 theme colors need to change.)

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 5 ---------------------------------------
$slide5 = $p.Slides.Item(5)
$tableShape = $slide5.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{202BF87D-3801-40E1-A598-9F15E327B125}")

# --- 2. Theme colors: Integral -> Office Theme ------------------------
# Theme color slot order exposed by ThemeColorScheme:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
# RGB() packs as R + G*256 + B*65536 (VBA/COM long color order).
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = $officeColors[$i - 1]
}
